$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7631788759517926
$ws.Range("B3").Value = 0.8588149599066598
$ws.Range("B4").Value = 0.754702226045021
$ws.Range("B5").Value = 0.8726784637810832
$ws.Range("B6").Value = 0.6643155405164123
$ws.Range("B7").Value = 0.9314458527886474
$ws.Range("B8").Value = 0.818754419330456
$ws.Range("B9").Value = 0.9691475873503521
$ws.Range("B10").Value = 0.5019974036702661

$ws.Range("A11").Value = "6_1"
$ws.Range("B11").Value = 0.9528021183690329
$ws.Range("A12").Value = "6_2"
$ws.Range("B12").Value = 0.2713279678068412
